$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-11-05 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-11-06 Thursday", 2) | Out-Null
$d.Content.Find.Execute("23+49=72", $true, $false, $false, $false, $false, $true, 1, $false, "89-32=57", 2) | Out-Null
$d.Content.Find.Execute("24-5=19", $true, $false, $false, $false, $false, $true, 1, $false, "74-73=1", 2) | Out-Null
$d.Content.Find.Execute("15+12=27", $true, $false, $false, $false, $false, $true, 1, $false, "95-37=58", 2) | Out-Null
$d.Content.Find.Execute("89-22=67", $true, $false, $false, $false, $false, $true, 1, $false, "95-25=70", 2) | Out-Null
$d.Content.Find.Execute("96-8=88", $true, $false, $false, $false, $false, $true, 1, $false, "19+71=90", 2) | Out-Null
$d.Content.Find.Execute("98-13=85", $true, $false, $false, $false, $false, $true, 1, $false, "20+26=46", 2) | Out-Null
$d.Content.Find.Execute("88-73=15", $true, $false, $false, $false, $false, $true, 1, $false, "6+20=26", 2) | Out-Null
$d.Content.Find.Execute("32+9=41", $true, $false, $false, $false, $false, $true, 1, $false, "7+32=39", 2) | Out-Null
$d.Content.Find.Execute("51+20=71", $true, $false, $false, $false, $false, $true, 1, $false, "78-19=59", 2) | Out-Null
$d.Content.Find.Execute("55+15=70", $true, $false, $false, $false, $false, $true, 1, $false, "89-15=74", 2) | Out-Null
$d.Content.Find.Execute("34+12=46", $true, $false, $false, $false, $false, $true, 1, $false, "25-7=18", 2) | Out-Null
$d.Content.Find.Execute("57+7=64", $true, $false, $false, $false, $false, $true, 1, $false, "25+53=78", 2) | Out-Null
$d.Content.Find.Execute("29-10=19", $true, $false, $false, $false, $false, $true, 1, $false, "20+24=44", 2) | Out-Null
$d.Content.Find.Execute("50-0=50", $true, $false, $false, $false, $false, $true, 1, $false, "50-15=35", 2) | Out-Null
$d.Content.Find.Execute("15+38=53", $true, $false, $false, $false, $false, $true, 1, $false, "41-31=10", 2) | Out-Null
$d.Content.Find.Execute("93+5=98", $true, $false, $false, $false, $false, $true, 1, $false, "63-2=61", 2) | Out-Null
$d.Content.Find.Execute("39-27=12", $true, $false, $false, $false, $false, $true, 1, $false, "25+41=66", 2) | Out-Null
$d.Content.Find.Execute("93-91=2", $true, $false, $false, $false, $false, $true, 1, $false, "98-30=68", 2) | Out-Null
$d.Content.Find.Execute("30-26=4", $true, $false, $false, $false, $false, $true, 1, $false, "72-25=47", 2) | Out-Null
$d.Content.Find.Execute("90-24=66", $true, $false, $false, $false, $false, $true, 1, $false, "75-10=65", 2) | Out-Null
$d.Content.Find.Execute("76-1=75", $true, $false, $false, $false, $false, $true, 1, $false, "73-13=60", 2) | Out-Null
$d.Content.Find.Execute("0+74=74", $true, $false, $false, $false, $false, $true, 1, $false, "54-8=46", 2) | Out-Null
$d.Content.Find.Execute("41+22=63", $true, $false, $false, $false, $false, $true, 1, $false, "33-0=33", 2) | Out-Null
$d.Content.Find.Execute("40-2=38", $true, $false, $false, $false, $false, $true, 1, $false, "70+10=80", 2) | Out-Null
$d.Content.Find.Execute("12+5=17", $true, $false, $false, $false, $false, $true, 1, $false, "79-28=51", 2) | Out-Null
$d.Content.Find.Execute("55+39=94", $true, $false, $false, $false, $false, $true, 1, $false, "13+21=34", 2) | Out-Null
$d.Content.Find.Execute("61+7=68", $true, $false, $false, $false, $false, $true, 1, $false, "30-17=13", 2) | Out-Null
$d.Content.Find.Execute("29+66=95", $true, $false, $false, $false, $false, $true, 1, $false, "68+22=90", 2) | Out-Null
$d.Content.Find.Execute("33+30=63", $true, $false, $false, $false, $false, $true, 1, $false, "86-3=83", 2) | Out-Null
$d.Content.Find.Execute("36+43=79", $true, $false, $false, $false, $false, $true, 1, $false, "0+61=61", 2) | Out-Null
$d.Content.Find.Execute("12+75=87", $true, $false, $false, $false, $false, $true, 1, $false, "17+50=67", 2) | Out-Null
$d.Content.Find.Execute("9+32=41", $true, $false, $false, $false, $false, $true, 1, $false, "87-45=42", 2) | Out-Null
$d.Content.Find.Execute("96-34=62", $true, $false, $false, $false, $false, $true, 1, $false, "93-8=85", 2) | Out-Null
$d.Content.Find.Execute("37-23=14", $true, $false, $false, $false, $false, $true, 1, $false, "55+9=64", 2) | Out-Null
$d.Content.Find.Execute("35+64=99", $true, $false, $false, $false, $false, $true, 1, $false, "85-21=64", 2) | Out-Null
$d.Content.Find.Execute("50-40=10", $true, $false, $false, $false, $false, $true, 1, $false, "55+27=82", 2) | Out-Null
$d.Content.Find.Execute("90-12=78", $true, $false, $false, $false, $false, $true, 1, $false, "29+36=65", 2) | Out-Null
$d.Content.Find.Execute("47-34=13", $true, $false, $false, $false, $false, $true, 1, $false, "79-72=7", 2) | Out-Null
$d.Content.Find.Execute("32-32=0", $true, $false, $false, $false, $false, $true, 1, $false, "82-63=19", 2) | Out-Null
$d.Content.Find.Execute("29+21=50", $true, $false, $false, $false, $false, $true, 1, $false, "45+8=53", 2) | Out-Null
$d.Content.Find.Execute("24-12=12", $true, $false, $false, $false, $false, $true, 1, $false, "7+81=88", 2) | Out-Null
$d.Content.Find.Execute("82+13=95", $true, $false, $false, $false, $false, $true, 1, $false, "7-4=3", 2) | Out-Null
$d.Content.Find.Execute("39+7=46", $true, $false, $false, $false, $false, $true, 1, $false, "91-34=57", 2) | Out-Null
$d.Content.Find.Execute("44-23=21", $true, $false, $false, $false, $false, $true, 1, $false, "53-43=10", 2) | Out-Null
$d.Content.Find.Execute("7+2=9", $true, $false, $false, $false, $false, $true, 1, $false, "49-23=26", 2) | Out-Null
$d.Content.Find.Execute("71-50=21", $true, $false, $false, $false, $false, $true, 1, $false, "59+22=81", 2) | Out-Null
$d.Content.Find.Execute("16+51=67", $true, $false, $false, $false, $false, $true, 1, $false, "47-13=34", 2) | Out-Null
$d.Content.Find.Execute("75-28=47", $true, $false, $false, $false, $false, $true, 1, $false, "92-23=69", 2) | Out-Null
$d.Content.Find.Execute("3-2=1", $true, $false, $false, $false, $false, $true, 1, $false, "62-38=24", 2) | Out-Null
$d.Content.Find.Execute("81-0=81", $true, $false, $false, $false, $false, $true, 1, $false, "58+27=85", 2) | Out-Null
$d.Content.Find.Execute("57-52=5", $true, $false, $false, $false, $false, $true, 1, $false, "88-70=18", 2) | Out-Null
$d.Content.Find.Execute("85-27=58", $true, $false, $false, $false, $false, $true, 1, $false, "86-1=85", 2) | Out-Null
$d.Content.Find.Execute("15+39=54", $true, $false, $false, $false, $false, $true, 1, $false, "81+5=86", 2) | Out-Null
$d.Content.Find.Execute("94-67=27", $true, $false, $false, $false, $false, $true, 1, $false, "56-41=15", 2) | Out-Null
$d.Content.Find.Execute("51-0=51", $true, $false, $false, $false, $false, $true, 1, $false, "83-65=18", 2) | Out-Null
$d.Content.Find.Execute("94-50=44", $true, $false, $false, $false, $false, $true, 1, $false, "61-6=55", 2) | Out-Null
$d.Content.Find.Execute("27+72=99", $true, $false, $false, $false, $false, $true, 1, $false, "75-44=31", 2) | Out-Null
$d.Content.Find.Execute("22+38=60", $true, $false, $false, $false, $false, $true, 1, $false, "40-16=24", 2) | Out-Null
$d.Content.Find.Execute("50-36=14", $true, $false, $false, $false, $false, $true, 1, $false, "58-5=53", 2) | Out-Null
$d.Content.Find.Execute("6+43=49", $true, $false, $false, $false, $false, $true, 1, $false, "85+11=96", 2) | Out-Null
$d.Content.Find.Execute("31-15=16", $true, $false, $false, $false, $false, $true, 1, $false, "99-55=44", 2) | Out-Null
$d.Content.Find.Execute("20+48=68", $true, $false, $false, $false, $false, $true, 1, $false, "42+20=62", 2) | Out-Null
$d.Content.Find.Execute("98-25=73", $true, $false, $false, $false, $false, $true, 1, $false, "13+49=62", 2) | Out-Null
$d.Content.Find.Execute("87-66=21", $true, $false, $false, $false, $false, $true, 1, $false, "2+66=68", 2) | Out-Null
$d.Content.Find.Execute("79-23=56", $true, $false, $false, $false, $false, $true, 1, $false, "41+11=52", 2) | Out-Null
$d.Content.Find.Execute("37+2=39", $true, $false, $false, $false, $false, $true, 1, $false, "10+74=84", 2) | Out-Null
$d.Content.Find.Execute("30+3=33", $true, $false, $false, $false, $false, $true, 1, $false, "86-20=66", 2) | Out-Null
$d.Content.Find.Execute("47+22=69", $true, $false, $false, $false, $false, $true, 1, $false, "26+11=37", 2) | Out-Null
$d.Content.Find.Execute("88-76=12", $true, $false, $false, $false, $false, $true, 1, $false, "16+12=28", 2) | Out-Null
$d.Content.Find.Execute("75-58=17", $true, $false, $false, $false, $false, $true, 1, $false, "79+10=89", 2) | Out-Null
$d.Content.Find.Execute("24+52=76", $true, $false, $false, $false, $false, $true, 1, $false, "95+3=98", 2) | Out-Null
$d.Content.Find.Execute("88-33=55", $true, $false, $false, $false, $false, $true, 1, $false, "33+22=55", 2) | Out-Null
$d.Content.Find.Execute("88-9=79", $true, $false, $false, $false, $false, $true, 1, $false, "28+28=56", 2) | Out-Null
$d.Content.Find.Execute("95-91=4", $true, $false, $false, $false, $false, $true, 1, $false, "37+32=69", 2) | Out-Null
$d.Content.Find.Execute("50-39=11", $true, $false, $false, $false, $false, $true, 1, $false, "99-55=44", 2) | Out-Null
$d.Content.Find.Execute("26+7=33", $true, $false, $false, $false, $false, $true, 1, $false, "71+0=71", 2) | Out-Null
$d.Content.Find.Execute("50+14=64", $true, $false, $false, $false, $false, $true, 1, $false, "62-32=30", 2) | Out-Null
$d.Content.Find.Execute("34+34=68", $true, $false, $false, $false, $false, $true, 1, $false, "25+45=70", 2) | Out-Null
$d.Content.Find.Execute("3+36=39", $true, $false, $false, $false, $false, $true, 1, $false, "1+47=48", 2) | Out-Null
$d.Content.Find.Execute("83-46=37", $true, $false, $false, $false, $false, $true, 1, $false, "8-4=4", 2) | Out-Null
$d.Content.Find.Execute("65+9=74", $true, $false, $false, $false, $false, $true, 1, $false, "61-61=0", 2) | Out-Null
$d.Content.Find.Execute("23+12=35", $true, $false, $false, $false, $false, $true, 1, $false, "61+27=88", 2) | Out-Null
$d.Content.Find.Execute("59+6=65", $true, $false, $false, $false, $false, $true, 1, $false, "84-14=70", 2) | Out-Null
$d.Content.Find.Execute("24+60=84", $true, $false, $false, $false, $false, $true, 1, $false, "42+34=76", 2) | Out-Null
$d.Content.Find.Execute("41-37=4", $true, $false, $false, $false, $false, $true, 1, $false, "18+27=45", 2) | Out-Null
$d.Content.Find.Execute("20+7=27", $true, $false, $false, $false, $false, $true, 1, $false, "77+11=88", 2) | Out-Null
$d.Content.Find.Execute("99-86=13", $true, $false, $false, $false, $false, $true, 1, $false, "62-3=59", 2) | Out-Null
$d.Content.Find.Execute("16+64=80", $true, $false, $false, $false, $false, $true, 1, $false, "7+48=55", 2) | Out-Null
$d.Content.Find.Execute("53-51=2", $true, $false, $false, $false, $false, $true, 1, $false, "27+47=74", 2) | Out-Null
$d.Content.Find.Execute("52-16=36", $true, $false, $false, $false, $false, $true, 1, $false, "36+52=88", 2) | Out-Null
$d.Content.Find.Execute("93-40=53", $true, $false, $false, $false, $false, $true, 1, $false, "54-42=12", 2) | Out-Null
$d.Content.Find.Execute("86-25=61", $true, $false, $false, $false, $false, $true, 1, $false, "88-65=23", 2) | Out-Null
$d.Content.Find.Execute("29+15=44", $true, $false, $false, $false, $false, $true, 1, $false, "94-34=60", 2) | Out-Null
$d.Content.Find.Execute("40+10=50", $true, $false, $false, $false, $false, $true, 1, $false, "42-19=23", 2) | Out-Null
$d.Content.Find.Execute("49-4=45", $true, $false, $false, $false, $false, $true, 1, $false, "61+31=92", 2) | Out-Null
$d.Content.Find.Execute("76-7=69", $true, $false, $false, $false, $false, $true, 1, $false, "78+17=95", 2) | Out-Null
$d.Content.Find.Execute("95-69=26", $true, $false, $false, $false, $false, $true, 1, $false, "83-18=65", 2) | Out-Null
$d.Content.Find.Execute("21+21=42", $true, $false, $false, $false, $false, $true, 1, $false, "89-9=80", 2) | Out-Null
$d.Content.Find.Execute("78-61=17", $true, $false, $false, $false, $false, $true, 1, $false, "12+69=81", 2) | Out-Null
$d.Content.Find.Execute("24-14=10", $true, $false, $false, $false, $false, $true, 1, $false, "31+24=55", 2) | Out-Null
